# Adding new scenarios in withdrawl test
#
# - Insert a new worksheet "withdrawlAmtMoreThanBalTest" right after
#   "withdrawlAmtLessThanBalTest" (and before "test_suite").
# - Populate it with a withdrawlAmount / errorMsg header row and one data row.
# - Update the selection on "withdrawlAmtLessThanBalTest" to the used range.
# - Leave the new sheet as the tab that is selected/active, matching the
#   authored workbook.

$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "withdrawlAmtLessThanBalTest" sheet ---
$lessSheet = $wb.Worksheets.Item("withdrawlAmtLessThanBalTest")
$lessSheet.Range("A1:A2").Select() | Out-Null

# --- Insert the new sheet right after it ---
$newSheet = $wb.Worksheets.Add($null, $lessSheet)
$newSheet.Name = "withdrawlAmtMoreThanBalTest"

# --- Fill in header + data row ---
$newSheet.Range("A1").Value = "withdrawlAmount"
$newSheet.Range("B1").Value = "errorMsg"
$newSheet.Range("A2").Value = 1067
$newSheet.Range("B2").Value = "Transaction Failed. You can not withdraw amount more than the balance."

# --- Size the columns to fit the new content ---
$newSheet.Columns.Item(1).ColumnWidth = 17.42578125
$newSheet.Columns.Item(2).ColumnWidth = 67.140625

# --- Leave the selection/active cell on the new sheet (it becomes the tab shown) ---
$newSheet.Range("I16").Select() | Out-Null
